# Auto commit at 2025-09-26 12:02:07.02
# Updates the Metrics sheet's raw figures; the "today" sheet pulls these in
# via =Metrics!Bn formulas (and E/F columns chain off of those), so they
# recalculate automatically once the source cells change.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value  = 376014.83999999997
$metrics.Range("B3").Value  = 303509.19000000006
$metrics.Range("B4").Value  = 118749.03
$metrics.Range("B5").Value  = 14909
$metrics.Range("B6").Value  = 4295265.7199999988
$metrics.Range("B7").Value  = 3631036.6699999995
$metrics.Range("B8").Value  = 1248114.7100000002
$metrics.Range("B9").Value  = 166069
$metrics.Range("B10").Value = 32760589.520999826
$metrics.Range("B11").Value = 19660906.740000006
$metrics.Range("B12").Value = 11529823.6
$metrics.Range("B13").Value = 1263696

# Update the cursor position on the Metrics sheet view.
$metrics.Range("F30").Select()

# Update the cursor position on the "today" sheet view; selecting a range
# here last also keeps it the active/visible tab, matching the workbook.
$today = $wb.Worksheets.Item("today")
$today.Range("H11").Select()
